$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Set-ParaXml($para, [string]$innerXml) {
    $full = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="' + $W + '"><w:body>' + $innerXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $start = $para.Range.Start
    $end = $para.Range.End - 1
    $r = $d.Range($start, $end)
    $r.InsertXML($full)
}

function Find-ParaWithText([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($needle)) {
            return $d.Paragraphs($i)
        }
    }
    return $null
}

$RFONT = '<w:rFonts w:ascii="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>'

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits after "2nd" in the
#    "2nd year student - CBCS" heading paragraph.
# ---------------------------------------------------------------------------
$pYear = Find-ParaWithText("year student - CBCS")
$innerYear = '<w:p><w:pPr><w:pStyle w:val="Heading3"/><w:contextualSpacing w:val="0"/></w:pPr>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>2nd</w:t></w:r>' + `
    '<w:r w:rsidR="00B371C1"><w:rPr>' + $RFONT + '</w:rPr><w:t xml:space="preserve"> year student - CBCS </w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pYear $innerYear

# ---------------------------------------------------------------------------
# 2) "Phone Number - 09897908973." paragraph: drop the gramStart/gramEnd
#    proofErr markers and split the number into "+91-" and the digits.
# ---------------------------------------------------------------------------
$pPhone = Find-ParaWithText("Phone Number - 09897908973")
$innerPhone = '<w:p>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>Phone Number - +91-</w:t></w:r>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>9897908973</w:t></w:r>' + `
    '<w:r w:rsidR="00B371C1"><w:rPr>' + $RFONT + '</w:rPr><w:t>.</w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pPhone $innerPhone

# ---------------------------------------------------------------------------
# 3) "Have a self coded portfolio blog.(link:www.chinmaychamoli.in)" paragraph:
#    split "Have a " into two runs, split "ling:www..." into "link" + the
#    relocated _GoBack bookmark + ":www.chinmaychamoli.in".
# ---------------------------------------------------------------------------
$pBlog = Find-ParaWithText("chinmaychamoli.in")
$innerBlog = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr><w:rPr>' + $RFONT + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>Have a</w:t></w:r>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>self coded</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t xml:space="preserve"> portfolio blog.(</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>link</w:t></w:r>' + `
    '<w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>:www.chinmaychamoli.in</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t xml:space="preserve">)  </w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pBlog $innerBlog

# ---------------------------------------------------------------------------
# 4) "... (Half Marathon) in years 2014." -> "... in year 2014." (typo fix)
# ---------------------------------------------------------------------------
$okYears = $d.Content.Find.Execute(" (Half Marathon) in years", $true, $false, $false, $false, $false, $true, 1, $false, " (Half Marathon) in year", 2)

# ---------------------------------------------------------------------------
# 5) "Have represented school...competitions and won medals in the same." ->
#    split the run into two at "competit" / "ions...".
# ---------------------------------------------------------------------------
$pRepresented = Find-ParaWithText("Have represented school")
$innerRepresented = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr>' + $RFONT + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>Have represented school in inter school badminton and chess competit</w:t></w:r>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>ions and won medals in the same.</w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pRepresented $innerRepresented

# ---------------------------------------------------------------------------
# 6) Delete the "Have taken part in many music and dance competitions in
#    school" paragraph entirely (everything below shifts up by one slot).
# ---------------------------------------------------------------------------
$pMusic = Find-ParaWithText("Have taken part in many music")
$pMusic.Range.Delete()

# ---------------------------------------------------------------------------
# 7) Move w:lastRenderedPageBreak from the "Coding." run to the
#    "Travelling & Hiking." run.
# ---------------------------------------------------------------------------
$pCoding = Find-ParaWithText("Coding.")
$innerCoding = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr>' + $RFONT + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:t>Coding.</w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pCoding $innerCoding

$pTravel = Find-ParaWithText("Travelling")
$innerTravel = '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:ind w:hanging="360"/><w:contextualSpacing/><w:rPr>' + $RFONT + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $RFONT + '</w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Travelling &amp; Hiking.</w:t></w:r>' + `
    '</w:p>'
Set-ParaXml $pTravel $innerTravel

Write-Output "Done"
